$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.123.00'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '2.380.82'
$ws.Range('E3').Value = '  +2.88%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.96'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.28'
$ws.Range('E6').Value = '  -4.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.571'
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('E9').Value = '  -2.64%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.68'
$ws.Range('E10').Value = '  -7.31%  '
$ws.Range('E11').Value = '  -1.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.17'
$ws.Range('E12').Value = '  -3.82%  '
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = '2.742.09'
$ws.Range('E14').Value = '  +2.88%  '
$ws.Range('D15').Value = '2.372.79'
$ws.Range('E15').Value = '  +2.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.819'
$ws.Range('E16').Value = '  -1.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.70'
$ws.Range('E17').Value = '  -2.64%  '
$ws.Range('D18').Value = '46.033.00'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.95'
$ws.Range('E19').Value = '  -3.09%  '
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.07'
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.61'
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '245.77'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('E24').Value = '  -4.68%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.95'
$ws.Range('E25').Value = '  -2.41%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '40.05'
$ws.Range('E27').Value = '  -8.05%  '
$ws.Range('E28').Value = '  -3.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.83'
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.81'
$ws.Range('E30').Value = '  +21.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '21.14'
$ws.Range('E31').Value = '  +4.25%  '
$ws.Range('E32').Value = '  +6.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.55'
$ws.Range('E33').Value = '  -4.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '147.07'
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0780'
$ws.Range('E35').Value = '  -3.11%  '
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.91'
$ws.Range('E37').Value = '  +5.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.07'
$ws.Range('E39').Value = '  -5.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.96'
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('E41').Value = '  -2.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.23'
$ws.Range('E42').Value = '  -7.09%  '
$ws.Range('D43').Value = '1.921.85'
$ws.Range('E43').Value = '  +4.02%  '
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.46'
$ws.Range('E45').Value = '  +3.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.80'
$ws.Range('E46').Value = '  -11.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.42'
$ws.Range('E47').Value = '  +5.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.188'
$ws.Range('E48').Value = '  -5.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '98.44'
$ws.Range('D50').Value = '2.612.98'
$ws.Range('E50').Value = '  +2.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '69.23'
$ws.Range('E51').Value = '  -8.27%  '
